$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows (2 and 3) have their Fecha (D), Volumen (J), Precio minimo (K),
# Precio maximo (L), Precio promedio ponderado (M) and Precio $/Kg (P) swapped.
$columns = @("D", "J", "K", "L", "M", "P")

foreach ($col in $columns) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $value2 = $cell2.Value2
    $value3 = $cell3.Value2

    $cell2.Value = $value3
    $cell3.Value = $value2
}
